$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-numeric-looking price cells to remain Text, matching the
# original inlineStr storage (avoids Excel auto-converting to a number).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.657.93'
$ws.Range("E2").Value = '  -0.51%  '
$ws.Range("D3").Value = '1.597.45'
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '211.16'
$ws.Range("E5").Value = '  +0.30%  '
$ws.Range("E6").Value = '  +1.18%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  -1.19%  '
$ws.Range("D10").Value = '19.67'
$ws.Range("E10").Value = '  +0.17%  '
$ws.Range("E11").Value = '  -0.23%  '
$ws.Range("E12").Value = '  -0.03%  '
$ws.Range("D13").Value = '1.611.29'
$ws.Range("E13").Value = '  +0.95%  '
$ws.Range("E14").Value = '  -0.59%  '
$ws.Range("D15").Value = '0.521'
$ws.Range("E15").Value = '  -1.45%  '
$ws.Range("D16").Value = '64.79'
$ws.Range("E16").Value = '  +2.09%  '
$ws.Range("D17").Value = '26.652.74'
$ws.Range("E17").Value = '  -0.37%  '
$ws.Range("E18").Value = '  -0.17%  '
$ws.Range("D19").Value = '209.63'
$ws.Range("E19").Value = '  +0.24%  '
$ws.Range("E20").Value = '  -0.05%  '
$ws.Range("E21").Value = '  +0.76%  '
$ws.Range("D22").Value = '4.26'
$ws.Range("E22").Value = '  -0.30%  '
$ws.Range("E23").Value = '  -1.49%  '
$ws.Range("D24").Value = '8.91'
$ws.Range("E24").Value = '  +0.74%  '
$ws.Range("D25").Value = '146.13'
$ws.Range("E25").Value = '  -0.16%  '
$ws.Range("E26").Value = '  +0.11%  '
$ws.Range("D27").Value = '7.18'
$ws.Range("E27").Value = '  -4.10%  '
$ws.Range("E28").Value = '  +2.42%  '
$ws.Range("D29").Value = '15.28'
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("E30").Value = '  +0.62%  '
$ws.Range("E31").Value = '  +0.17%  '
$ws.Range("E32").Value = '  -0.72%  '
$ws.Range("E33").Value = '  -0.68%  '
$ws.Range("E34").Value = '  -0.80%  '
$ws.Range("D35").Value = '1.295.40'
$ws.Range("E35").Value = '  -1.45%  '
$ws.Range("E36").Value = '  +0.82%  '
$ws.Range("D37").Value = '1.48'
$ws.Range("E37").Value = '  -2.23%  '
$ws.Range("E38").Value = '  -1.09%  '
$ws.Range("E39").Value = '  +2.50%  '
$ws.Range("E41").Value = '  +2.06%  '
$ws.Range("E42").Value = '  +1.01%  '
$ws.Range("E43").Value = '  -0.08%  '
$ws.Range("D44").Value = '63.80'
$ws.Range("E44").Value = '  +1.65%  '
$ws.Range("D45").Value = '1.733.93'
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("D46").Value = '0.894'
$ws.Range("E46").Value = '  +10.59%  '
$ws.Range("D47").Value = '90.19'
$ws.Range("E47").Value = '  +1.30%  '
$ws.Range("E48").Value = '  +0.83%  '
$ws.Range("E49").Value = '  +2.63%  '
$ws.Range("E50").Value = '  -1.16%  '
$ws.Range("D51").Value = '7.48'
$ws.Range("E51").Value = '  +1.03%  '
